# Update the StructureDefinition-reference workbook:
#  - bump Version 5.0.0 -> 6.0.0
#  - bump Date to the new publish timestamp
#  - change Publisher from "Contact"/"No display for ContactDetail" to
#    "Alvearie Team", and turn the (duplicate) second Contact row into a
#    single new "Jurisdiction" / "United States of America" row
#  - update the root Extension element's Short/Definition on the Elements
#    sheet to match the new Title/Description ("Insight Reference" /
#    "Reference to content leveraged to produce the insight.")

$wb = $excel.ActiveWorkbook

# --- Metadata sheet (Property/Value table) ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The old sheet had two identical "Contact" rows (10 & 11); row 10 was
# repurposed above into the new "Jurisdiction" row, so the now-redundant
# duplicate row 11 is removed, shifting everything below it up by one.
$meta.Rows.Item(11).Delete()

# --- Elements sheet (element definitions table) ---
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Insight Reference"
$elements.Range("L2").Value = "Reference to content leveraged to produce the insight."
